$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 23; shifts existing rows 23..67 down to 24..68
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new weekly record
$ws.Cells.Item(23, 1).Value = 3
$ws.Cells.Item(23, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value = 44526
$ws.Cells.Item(23, 5).Value = 5
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103003
$ws.Cells.Item(23, 10).Value = "Damasco"
$ws.Cells.Item(23, 11).Value = "Castle Brite"
$ws.Cells.Item(23, 12).Value = "Segunda"
$ws.Cells.Item(23, 13).Value = 60
$ws.Cells.Item(23, 14).Value = 18000
$ws.Cells.Item(23, 15).Value = 18000
$ws.Cells.Item(23, 16).Value = 18000
$ws.Cells.Item(23, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(23, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(23, 19).Value = 1000
$ws.Cells.Item(23, 20).Value = 18
